$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.614.36'
$ws.Range("E2").Value = '  +2.66%  '
$ws.Range("D3").Value = '2.652.65'
$ws.Range("E3").Value = '  +2.14%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = "'599.34"
$ws.Range("E5").Value = '  +1.74%  '
$ws.Range("D6").Value = "'154.61"
$ws.Range("E6").Value = '  +3.40%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  +0.32%  '
$ws.Range("D9").Value = '2.651.15'
$ws.Range("E9").Value = '  +2.15%  '
$ws.Range("E10").Value = '  +10.64%  '
$ws.Range("E11").Value = '  -0.41%  '
$ws.Range("E12").Value = '  +1.37%  '
$ws.Range("E13").Value = '  +1.44%  '
$ws.Range("D14").Value = "'27.92"
$ws.Range("E14").Value = '  +2.77%  '
$ws.Range("E15").Value = '  +5.98%  '
$ws.Range("D16").Value = '3.133.81'
$ws.Range("E16").Value = '  +2.10%  '
$ws.Range("D17").Value = '68.480.90'
$ws.Range("E17").Value = '  +2.56%  '
$ws.Range("D18").Value = '2.647.65'
$ws.Range("E18").Value = '  +1.94%  '
$ws.Range("E19").Value = '  +3.80%  '
$ws.Range("D20").Value = "'368.19"
$ws.Range("E20").Value = '  +1.72%  '
$ws.Range("D21").Value = "'7.44"
$ws.Range("E21").Value = '  +1.85%  '
$ws.Range("E22").Value = '  -0.21%  '
$ws.Range("E23").Value = '  +0.84%  '
$ws.Range("E24").Value = '  +4.70%  '
$ws.Range("D25").Value = "'72.67"
$ws.Range("E25").Value = '  +0.41%  '
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("D27").Value = "'9.99"
$ws.Range("E27").Value = '  +0.86%  '
$ws.Range("D28").Value = "'0.0000105"
$ws.Range("E28").Value = '  +7.70%  '
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = '  +0.05%  '
$ws.Range("D31").Value = "'574.34"
$ws.Range("E31").Value = '  -1.25%  '
$ws.Range("E32").Value = '  +4.76%  '
$ws.Range("D33").Value = "'7.98"
$ws.Range("E33").Value = '  +5.08%  '
$ws.Range("E34").Value = '  +2.90%  '
$ws.Range("E35").Value = '  +4.73%  '
$ws.Range("E36").Value = '  +0.03%  '
$ws.Range("E37").Value = '  +3.98%  '
$ws.Range("D38").Value = "'158.31"
$ws.Range("E38").Value = '  +1.34%  '
$ws.Range("E39").Value = '  +5.33%  '
$ws.Range("E40").Value = '  +1.73%  '
$ws.Range("E41").Value = '  +3.96%  '
$ws.Range("E42").Value = '  +0.83%  '
$ws.Range("E43").Value = '  +6.41%  '
$ws.Range("E44").Value = '  +4.59%  '
$ws.Range("D45").Value = '0.0₆0321'
$ws.Range("E45").Value = '  +12.78%  '
$ws.Range("E46").Value = '  +0.09%  '
$ws.Range("D47").Value = "'40.58"
$ws.Range("E47").Value = '  -0.42%  '
$ws.Range("D48").Value = "'156.89"
$ws.Range("E48").Value = '  +3.29%  '
$ws.Range("E49").Value = '  +0.86%  '
$ws.Range("E50").Value = '  +2.65%  '
$ws.Range("E51").Value = '  +3.25%  '
